$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.313393592834473
$ws.Range("B1").Value = 3.328510522842407
$ws.Range("C1").Value = 3.047177314758301
$ws.Range("D1").Value = 3.45962381362915
$ws.Range("E1").Value = 1.736697316169739
